$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated K (strikeouts) column values, regenerated from source data
$kValues = @{
    2 = 5
    3 = 4
    4 = 2
    5 = 4
    6 = 5
    7 = 3
    8 = 1
    9 = 4
    10 = 1
    11 = 2
    12 = 7
    13 = 5
    14 = 5
    15 = 2
    16 = 3
    17 = 0
    18 = 2
    19 = 5
    20 = 2
    21 = 2
    22 = 4
    23 = 5
    24 = 2
    25 = 3
    26 = 2
    27 = 5
    28 = 1
    29 = 2
    30 = 1
    31 = 6
    32 = 6
    33 = 1
    34 = 5
    35 = 7
    36 = 1
    37 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

